$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.685.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").Value = '''1.849.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("D4").Value = '''1.029'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''322.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '''0.4385'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '''0.3791'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.56%  '
$ws.Range("D9").Value = '''0.07383'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").Value = '''0.8814'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").Value = '''21.51'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '''1.852.72'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = '''5.500'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").Value = '''6.707'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").Value = '''0.07140'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").Value = '''85.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.90%  '
$ws.Range("D17").Value = '''1.035'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '''0.000009056'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '''15.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").Value = '''27.717.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").Value = '''5.290'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.96%  '
$ws.Range("D23").Value = '''11.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.08%  '
$ws.Range("D24").Value = '''2.084.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").Value = '''2.040'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.05%  '
$ws.Range("D26").Value = '''157.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.47%  '
$ws.Range("D27").Value = '''18.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("D28").Value = '''2.000'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.10%  '
$ws.Range("D29").Value = '''5.329'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.48%  '
$ws.Range("D30").Value = '''117.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("D31").Value = '''0.09036'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '''1.207'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7693'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("D34").Value = '''2.996'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.30%  '
$ws.Range("D35").Value = '''4.543'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("D36").Value = '''1.030'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").Value = '''1.148'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("D38").Value = '''0.01973'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").Value = '''0.05263'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("D40").Value = '''2.848'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("D41").Value = '''0.5177'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").Value = '''0.1671'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").Value = '''6.854'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.54%  '
$ws.Range("D44").Value = '''8.758'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.65%  '
$ws.Range("D45").Value = '''110.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.18%  '
$ws.Range("D46").Value = '''10.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").Value = '''0.06618'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.42%  '
$ws.Range("D48").Value = '''1.032'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").Value = '''1.700'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").Value = '''0.4692'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.01%  '
$ws.Range("D51").Value = '''1.894'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.17%  '
